# Add new "Plays  May 11, 2025" column (F) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: same text/number format as the other header cells (A1:E1)
$ws.Range("F1").Value = "Plays  May 11, 2025"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data rows F2:F91
$data = @{
    2 = '1841991'
    3 = '49152'
    4 = '123928'
    5 = '505808'
    6 = '102997'
    7 = '5690'
    8 = '66632'
    9 = '30722'
    10 = '15535'
    11 = '72591'
    12 = 'Invalid URL'
    13 = '57556'
    14 = '62694'
    15 = '32419'
    16 = '85842'
    17 = '15303'
    18 = '9810'
    19 = 'Invalid URL'
    20 = '98292'
    21 = '5879'
    22 = '5347'
    23 = '10386'
    24 = '7145'
    25 = '45358'
    26 = '28824'
    27 = '17721'
    28 = '20523'
    29 = '105671'
    30 = '3685'
    31 = '17515'
    32 = '5003'
    33 = '6544'
    34 = '16523'
    35 = '19360'
    36 = '6594'
    37 = '52795'
    38 = '6050'
    39 = '4372'
    40 = '9913'
    41 = '21007'
    42 = '13773'
    43 = '11633'
    44 = '29316'
    45 = '8016'
    46 = '14177'
    47 = 'Invalid URL'
    48 = '14092'
    49 = '26343'
    50 = '7259'
    51 = '3874'
    52 = '33175'
    53 = '7239'
    54 = '11638'
    55 = '17344'
    56 = '7351'
    57 = '15794'
    58 = '6086'
    59 = '23017'
    60 = '4929'
    61 = '10333'
    62 = '8828'
    63 = '9202'
    64 = '6446'
    65 = '9575'
    66 = '4672'
    67 = 'Invalid URL'
    68 = '12193'
    69 = '5077'
    70 = '6353'
    71 = 'Invalid URL'
    72 = '1629'
    73 = '45358'
    74 = '105671'
    75 = '1150'
    76 = '14177'
    77 = '3874'
    78 = '29316'
    79 = '1118'
    80 = '4229'
    81 = '5490'
    82 = '2125'
    83 = '2068'
    84 = '2120'
    85 = '1322'
    86 = '27487'
    87 = '416'
    88 = '1161'
    89 = '1061'
    90 = '504'
    91 = '1420'
}

foreach ($row in $data.Keys | Sort-Object) {
    $value = $data[$row]
    $cell = $ws.Cells.Item($row, 6)
    if ($value -match '^-?\d+$') {
        # Numeric-looking play counts are stored as text in this sheet,
        # matching the existing columns (A-E), so force text entry, then
        # drop the now-unneeded "Text" number format so the cell keeps the
        # default (unstyled) look of its neighbours in columns A-E.
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.ClearFormats()
    } else {
        $cell.Value = $value
    }
}

Write-Host ("Final used range: " + $ws.UsedRange.Address())